$wb = $excel.ActiveWorkbook

# Sheet 1: BWE_echoes_detection
$ws = $wb.Worksheets.Item(1)
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 100
$ws.Range("E4").Value = 100
$ws.Range("F4").Value = 100
$ws.Range("B5").Value = 100
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = 100
$ws.Range("F5").Value = 100
$ws.Range("B6").Value = 100
$ws.Range("D6").Value = 100
$ws.Range("E6").Value = 100
$ws.Range("F6").Value = 100
$ws.Range("B7").Value = 100
$ws.Range("E7").Value = 100
$ws.Range("F7").Value = 100
$ws.Range("B8").Value = 100
$ws.Range("D8").Value = 100

# Sheet 2: BWE_mean_distance_error
$ws = $wb.Worksheets.Item(2)
$ws.Range("B4").Value = 0.02365783452653334
$ws.Range("C4").Value = 0.001543287327478072
$ws.Range("D4").Value = 0.004641670972027444
$ws.Range("E4").Value = 0.005618864860875331
$ws.Range("F4").Value = -0.004364159716584307
$ws.Range("B5").Value = 0.02181784633552293
$ws.Range("D5").Value = 0.004694811425197448
$ws.Range("E5").Value = 0.005029891504908117
$ws.Range("F5").Value = -0.00288508377002
$ws.Range("B6").Value = 0.01946638128275146
$ws.Range("D6").Value = 0.004305114768617639
$ws.Range("E6").Value = 0.004161930769798526
$ws.Range("F6").Value = -0.001764705882353011
$ws.Range("B7").Value = 0.01666765074913267
$ws.Range("E7").Value = 0.00287327478042663
$ws.Range("F7").Value = -0.0006044726548085195
$ws.Range("B8").Value = 0.01181858439737245
$ws.Range("D8").Value = 0.003082884345708216

# Sheet 3: BWE_STD_distance_error
$ws = $wb.Worksheets.Item(3)
$ws.Range("B4").Value = 0.002035967391253168
$ws.Range("C4").Value = 0.002490319006590181
$ws.Range("D4").Value = 0.002043143839331984
$ws.Range("E4").Value = 0.001880883832005884
$ws.Range("F4").Value = 0.002546193042046593
$ws.Range("B5").Value = 0.001704259925763946
$ws.Range("D5").Value = 0.001279636429058974
$ws.Range("E5").Value = 0.001342780772921734
$ws.Range("F5").Value = 0.001592907613271836
$ws.Range("B6").Value = 0.001788633010541431
$ws.Range("D6").Value = 0.001548941093028976
$ws.Range("E6").Value = 0.001370562969376771
$ws.Range("F6").Value = 0.001451599173634971
$ws.Range("B7").Value = 0.001161869472433352
$ws.Range("E7").Value = 0.001225785198558041
$ws.Range("F7").Value = 0.001359645151408854
$ws.Range("B8").Value = 0.001249771922356871
$ws.Range("D8").Value = 0.00156170102840523

# Sheet 4: BWE_mean_amplitude_error_1
$ws = $wb.Worksheets.Item(4)
$ws.Range("B4").Value = -26.54589195190782
$ws.Range("C4").Value = -4.614615457333819
$ws.Range("D4").Value = -9.156226171518295
$ws.Range("E4").Value = -6.116429325715128
$ws.Range("F4").Value = -2.310349471226028
$ws.Range("B5").Value = -24.21109122750564
$ws.Range("D5").Value = -8.054217653056249
$ws.Range("E5").Value = -4.744366658482702
$ws.Range("F5").Value = -0.5939367701619458
$ws.Range("B6").Value = -22.30166806259818
$ws.Range("D6").Value = -7.602445697231851
$ws.Range("E6").Value = -4.452177570026453
$ws.Range("F6").Value = 0.7389267850343705
$ws.Range("B7").Value = -18.04829467257114
$ws.Range("E7").Value = 0.4039858422133276
$ws.Range("F7").Value = 2.224410226504793
$ws.Range("B8").Value = -10.99568223481967
$ws.Range("D8").Value = -2.149683673259656

# Sheet 5: BWE_STD_amplitude_error_1
$ws = $wb.Worksheets.Item(5)
$ws.Range("B4").Value = 3.334744370426787
$ws.Range("C4").Value = 5.299820319418025
$ws.Range("D4").Value = 3.737126619954745
$ws.Range("E4").Value = 3.396576810952055
$ws.Range("F4").Value = 5.016290421391401
$ws.Range("B5").Value = 2.526736109056072
$ws.Range("D5").Value = 2.519867400252476
$ws.Range("E5").Value = 2.289296289438998
$ws.Range("F5").Value = 3.303014865323734
$ws.Range("B6").Value = 2.847031525659731
$ws.Range("D6").Value = 2.969924482433771
$ws.Range("E6").Value = 2.462110949048284
$ws.Range("F6").Value = 3.410486569297024
$ws.Range("B7").Value = 1.568207772039397
$ws.Range("E7").Value = 3.226978287566372
$ws.Range("F7").Value = 2.006297947669716
$ws.Range("B8").Value = 1.365072671008647
$ws.Range("D8").Value = 1.166648896344525

# Sheet 6: BWE_mean_amplitude_error_2
$ws = $wb.Worksheets.Item(6)
$ws.Range("B4").Value = -26.86031049084387
$ws.Range("C4").Value = -4.647840414730961
$ws.Range("D4").Value = -9.00760275339627
$ws.Range("E4").Value = -6.325938980223003
$ws.Range("F4").Value = -3.230061096220323
$ws.Range("B5").Value = -23.51424401310565
$ws.Range("D5").Value = -8.250225357458735
$ws.Range("E5").Value = -4.994416566227185
$ws.Range("F5").Value = -1.659610840162146
$ws.Range("B6").Value = -20.86382596236752
$ws.Range("D6").Value = -8.536472800898901
$ws.Range("E6").Value = -4.179023357988734
$ws.Range("F6").Value = -0.021752606694643
$ws.Range("B7").Value = -14.16240271748426
$ws.Range("E7").Value = 0.1024409273591615
$ws.Range("F7").Value = 5.208189141637962
$ws.Range("B8").Value = -10.53113940143941
$ws.Range("D8").Value = -2.776433776219893

# Sheet 7: BWE_STD_amplitude_error_2
$ws = $wb.Worksheets.Item(7)
$ws.Range("B4").Value = 3.452821354779885
$ws.Range("C4").Value = 5.059483158385259
$ws.Range("D4").Value = 3.680676563626293
$ws.Range("E4").Value = 3.211778697013214
$ws.Range("F4").Value = 5.182582469648143
$ws.Range("B5").Value = 2.407673343220441
$ws.Range("D5").Value = 2.666157798265959
$ws.Range("E5").Value = 2.266231230746703
$ws.Range("F5").Value = 3.724977976817463
$ws.Range("B6").Value = 2.508002765111691
$ws.Range("D6").Value = 3.139726948013384
$ws.Range("E6").Value = 2.353839641458567
$ws.Range("F6").Value = 3.474467643659124
$ws.Range("B7").Value = 1.333850090809651
$ws.Range("E7").Value = 1.455138661921824
$ws.Range("F7").Value = 2.510831843449406
$ws.Range("B8").Value = 2.552597245652156
$ws.Range("D8").Value = 1.166648896344525
